# Daily "Model holdings" upload refresh (2021-07-08 -> 2021-07-09):
#  - bump the disclaimer date in the confidential notice (A38)
#  - refresh the Weight (D) / Percent Change (E) figures for each
#    holding (rows 2-35) to the new day's values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected (no known password); lift it for the
# duration of the refresh and restore protection afterwards.
$ws.Unprotect()

# New Weight / Percent Change values, keyed by row number.
$newValues = @(
    @(2, 0.09452737317642468, 0.01305501256632224),
    @(3, 0.07959819295862139, 0.001874414245548239),
    @(4, 0.05353137718760893, -0.003234702163525283),
    @(5, 0.05111964547787751, 0.02331265961327977),
    @(6, 0.04752220425765116, 0.007861037149740069),
    @(7, 0.04114275712967917, 0.03199946998807479),
    @(8, 0.03587800605694562, 0.003794664278174098),
    @(9, 0.03881035964369895, 0.003962621244381292),
    @(10, 0.03404385044996708, 0.005086324235260564),
    @(11, 0.03528723149382183, 0.005268935236004468),
    @(12, 0.03471123346785242, 0.01380008679299882),
    @(13, 0.03004256503469617, 0.03249097472924189),
    @(14, 0.03198876409330164, 0.01114459722483807),
    @(15, 0.03186481333503574, 0.007428617736837895),
    @(16, 0.02974818198381467, 0.0245370370370368),
    @(17, 0.02947703970010802, 0.0003650167907722768),
    @(18, 0.02809593559151339, 0.006341814616810604),
    @(19, 0.02422591747232268, -0.003434655675768394),
    @(20, 0.02076863816277431, 0.007004310344827402),
    @(21, 0.02156944039962642, 0.01812437645493858),
    @(22, 0.02174360269191207, 0.006096434509514248),
    @(23, 0.0206604681723432, 0.01083228019498095),
    @(24, 0.01859663066260575, 0.009581263307310239),
    @(25, 0.02210196958328731, 0.02535342913891814),
    @(26, 0.02027111359140149, 0.009171974522293125),
    @(27, 0.019865978242625, 0.009012392039053685),
    @(28, 0.01882014371975652, 0.004024819721616391),
    @(29, 0.0207398023960944, 0.002303430243416393),
    @(30, 0.01142111552813209, 0.007411036163344153),
    @(31, 0.008465664635729666, 0.01748856126080311),
    @(32, 0.007614363941806266, 0.009834953651367995),
    @(33, 0.008693051154076243, -0.002392936710949867),
    @(34, 0.007052568606888162, -0.003010577705451656),
    @(35, 0.9999999999999999, 0.01001031488370274)
)

foreach ($entry in $newValues) {
    $rowNum = $entry[0]
    $ws.Cells.Item($rowNum, 4).Value = $entry[1]
    $ws.Cells.Item($rowNum, 5).Value = $entry[2]
}

# Bump the "as of" date in the confidential disclaimer text.
$ws.Range("A38").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-09 for illustrative purposes only and are subject to change."

# Restore sheet protection.
$ws.Protect()
